$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the array-of-scores shared string (all Z2:Z17 cells reference the same text,
# so assigning the same new text to every cell in the range keeps them deduplicated
# onto a single shared-string entry, matching the original single <si> replacement.)
$newArrayStr = "[0.64384921 0.55753968 0.69032922 0.72016461 0.66326531]"
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("Z$r").Value = $newArrayStr
}

# Update numeric metric cells (B:E, R:Y, AA:AB) for rows 2-17
$ws.Range("B2").Value = 0.01667771339416504
$ws.Range("C2").Value = 0.004421853766433573
$ws.Range("D2").Value = 0.002419042587280274
$ws.Range("E2").Value = 0.0009763600723563534
$ws.Range("R2").Value = 0.5614035087719299
$ws.Range("S2").Value = 0.4313725490196078
$ws.Range("T2").Value = 0.5862068965517241
$ws.Range("U2").Value = 0.6415094339622641
$ws.Range("V2").Value = 0.5964912280701754
$ws.Range("W2").Value = 0.5633967232751402
$ws.Range("X2").Value = 0.07092460971388945
$ws.Range("Y2").Value = 5
$ws.Range("AA2").Value = 0.6550296044343662
$ws.Range("AB2").Value = 0.05511033388910008
$ws.Range("B3").Value = 0.006073236465454102
$ws.Range("C3").Value = 0.003264044651622484
$ws.Range("D3").Value = 0.001746082305908203
$ws.Range("E3").Value = 0.0005268508702025855
$ws.Range("R3").Value = 0.5333333333333333
$ws.Range("S3").Value = 0.4210526315789473
$ws.Range("T3").Value = 0.6268656716417911
$ws.Range("U3").Value = 0.5882352941176471
$ws.Range("V3").Value = 0.6363636363636364
$ws.Range("W3").Value = 0.561170113407071
$ws.Range("X3").Value = 0.0788800244827743
$ws.Range("Y3").Value = 7
$ws.Range("AA3").Value = 0.6550296044343662
$ws.Range("AB3").Value = 0.05511033388910008
$ws.Range("B4").Value = 0.001752185821533203
$ws.Range("C4").Value = 0.0002385684210073639
$ws.Range("D4").Value = 0.0009346485137939453
$ws.Range("E4").Value = 0.000273238692802794
$ws.Range("R4").Value = 0.5818181818181818
$ws.Range("S4").Value = 0.44
$ws.Range("T4").Value = 0.5614035087719299
$ws.Range("U4").Value = 0.6296296296296297
$ws.Range("V4").Value = 0.5185185185185186
$ws.Range("W4").Value = 0.546273967747652
$ws.Range("X4").Value = 0.0640390220350366
$ws.Range("Y4").Value = 11
$ws.Range("AA4").Value = 0.6550296044343662
$ws.Range("AB4").Value = 0.05511033388910008
$ws.Range("B5").Value = 0.00205235481262207
$ws.Range("C5").Value = 0.0002030336223793909
$ws.Range("D5").Value = 0.0008274078369140625
$ws.Range("E5").Value = 0.00004983817485585115
$ws.Range("R5").Value = 0.5
$ws.Range("S5").Value = 0.4363636363636363
$ws.Range("T5").Value = 0.6060606060606061
$ws.Range("U5").Value = 0.5599999999999999
$ws.Range("V5").Value = 0.5937499999999999
$ws.Range("W5").Value = 0.5392348484848485
$ws.Range("X5").Value = 0.06323074371958229
$ws.Range("Y5").Value = 13
$ws.Range("AA5").Value = 0.6550296044343662
$ws.Range("AB5").Value = 0.05511033388910008
$ws.Range("B6").Value = 0.002403783798217774
$ws.Range("C6").Value = 0.0007959683644120436
$ws.Range("D6").Value = 0.0009770870208740234
$ws.Range("E6").Value = 0.0003896888886730941
$ws.Range("R6").Value = 0.634920634920635
$ws.Range("S6").Value = 0.4745762711864406
$ws.Range("T6").Value = 0.6129032258064516
$ws.Range("U6").Value = 0.6545454545454545
$ws.Range("V6").Value = 0.6774193548387097
$ws.Range("W6").Value = 0.6108729882595384
$ws.Range("X6").Value = 0.07140715576988323
$ws.Range("Y6").Value = 1
$ws.Range("AA6").Value = 0.6550296044343662
$ws.Range("AB6").Value = 0.05511033388910008
$ws.Range("B7").Value = 0.002150058746337891
$ws.Range("C7").Value = 0.000135158649256228
$ws.Range("D7").Value = 0.0008048057556152344
$ws.Range("E7").Value = 0.00003689579701401173
$ws.Range("R7").Value = 0.5333333333333333
$ws.Range("S7").Value = 0.4210526315789473
$ws.Range("T7").Value = 0.6268656716417911
$ws.Range("U7").Value = 0.5882352941176471
$ws.Range("V7").Value = 0.6363636363636364
$ws.Range("W7").Value = 0.561170113407071
$ws.Range("X7").Value = 0.0788800244827743
$ws.Range("Y7").Value = 7
$ws.Range("AA7").Value = 0.6550296044343662
$ws.Range("AB7").Value = 0.05511033388910008
$ws.Range("B8").Value = 0.001673269271850586
$ws.Range("C8").Value = 0.0004173958431448494
$ws.Range("D8").Value = 0.0007188320159912109
$ws.Range("E8").Value = 0.0000477487597301953
$ws.Range("R8").Value = 0.6229508196721312
$ws.Range("S8").Value = 0.4363636363636363
$ws.Range("T8").Value = 0.6
$ws.Range("U8").Value = 0.6666666666666666
$ws.Range("V8").Value = 0.6000000000000001
$ws.Range("W8").Value = 0.5851962245404868
$ws.Range("X8").Value = 0.07829723524227525
$ws.Range("Y8").Value = 2
$ws.Range("AA8").Value = 0.6550296044343662
$ws.Range("AB8").Value = 0.05511033388910008
$ws.Range("B9").Value = 0.001863336563110351
$ws.Range("C9").Value = 0.0001607630578983279
$ws.Range("D9").Value = 0.00070648193359375
$ws.Range("E9").Value = 0.00002606925916629189
$ws.Range("R9").Value = 0.5
$ws.Range("S9").Value = 0.4363636363636363
$ws.Range("T9").Value = 0.6060606060606061
$ws.Range("U9").Value = 0.5599999999999999
$ws.Range("V9").Value = 0.5937499999999999
$ws.Range("W9").Value = 0.5392348484848485
$ws.Range("X9").Value = 0.06323074371958229
$ws.Range("Y9").Value = 13
$ws.Range("AA9").Value = 0.6550296044343662
$ws.Range("AB9").Value = 0.05511033388910008
$ws.Range("B10").Value = 0.001822853088378906
$ws.Range("C10").Value = 0.0001810898331769118
$ws.Range("D10").Value = 0.0007233142852783204
$ws.Range("E10").Value = 0.00002994613282036806
$ws.Range("R10").Value = 0.5937499999999999
$ws.Range("S10").Value = 0.4482758620689655
$ws.Range("T10").Value = 0.6363636363636364
$ws.Range("U10").Value = 0.6153846153846153
$ws.Range("V10").Value = 0.6060606060606061
$ws.Range("W10").Value = 0.5799669439755646
$ws.Range("X10").Value = 0.06730288440842029
$ws.Range("Y10").Value = 3
$ws.Range("AA10").Value = 0.6550296044343662
$ws.Range("AB10").Value = 0.05511033388910008
$ws.Range("B11").Value = 0.001897716522216797
$ws.Range("C11").Value = 0.0001443518946606092
$ws.Range("D11").Value = 0.0006984710693359375
$ws.Range("E11").Value = 0.00003377899081573861
$ws.Range("R11").Value = 0.5333333333333333
$ws.Range("S11").Value = 0.4210526315789473
$ws.Range("T11").Value = 0.6268656716417911
$ws.Range("U11").Value = 0.5882352941176471
$ws.Range("V11").Value = 0.6363636363636364
$ws.Range("W11").Value = 0.561170113407071
$ws.Range("X11").Value = 0.0788800244827743
$ws.Range("Y11").Value = 7
$ws.Range("AA11").Value = 0.6550296044343662
$ws.Range("AB11").Value = 0.05511033388910008
$ws.Range("B12").Value = 0.002337169647216797
$ws.Range("C12").Value = 0.0009706953475183825
$ws.Range("D12").Value = 0.001106977462768555
$ws.Range("E12").Value = 0.0007212996552790281
$ws.Range("R12").Value = 0.5901639344262296
$ws.Range("S12").Value = 0.4363636363636363
$ws.Range("T12").Value = 0.6129032258064516
$ws.Range("U12").Value = 0.5882352941176471
$ws.Range("V12").Value = 0.5806451612903226
$ws.Range("W12").Value = 0.5616622504008575
$ws.Range("X12").Value = 0.06356760777417875
$ws.Range("Y12").Value = 6
$ws.Range("AA12").Value = 0.6550296044343662
$ws.Range("AB12").Value = 0.05511033388910008
$ws.Range("B13").Value = 0.002263164520263672
$ws.Range("C13").Value = 0.00111588786048087
$ws.Range("D13").Value = 0.001127052307128906
$ws.Range("E13").Value = 0.0005560288776991021
$ws.Range("R13").Value = 0.5
$ws.Range("S13").Value = 0.4363636363636363
$ws.Range("T13").Value = 0.6060606060606061
$ws.Range("U13").Value = 0.5599999999999999
$ws.Range("V13").Value = 0.5937499999999999
$ws.Range("W13").Value = 0.5392348484848485
$ws.Range("X13").Value = 0.06323074371958229
$ws.Range("Y13").Value = 13
$ws.Range("AA13").Value = 0.6550296044343662
$ws.Range("AB13").Value = 0.05511033388910008
$ws.Range("B14").Value = 0.002850341796875
$ws.Range("C14").Value = 0.001238884314968684
$ws.Range("D14").Value = 0.0007994174957275391
$ws.Range("E14").Value = 0.0001131291510046847
$ws.Range("R14").Value = 0.5806451612903226
$ws.Range("S14").Value = 0.440677966101695
$ws.Range("T14").Value = 0.6176470588235294
$ws.Range("U14").Value = 0.5882352941176471
$ws.Range("V14").Value = 0.6268656716417911
$ws.Range("W14").Value = 0.570814230394997
$ws.Range("X14").Value = 0.06733597230629797
$ws.Range("Y14").Value = 4
$ws.Range("AA14").Value = 0.6550296044343662
$ws.Range("AB14").Value = 0.05511033388910008
$ws.Range("B15").Value = 0.002422428131103516
$ws.Range("C15").Value = 0.0008598976573005624
$ws.Range("D15").Value = 0.00081787109375
$ws.Range("E15").Value = 0.0002297692975159042
$ws.Range("R15").Value = 0.5333333333333333
$ws.Range("S15").Value = 0.4210526315789473
$ws.Range("T15").Value = 0.6268656716417911
$ws.Range("U15").Value = 0.5882352941176471
$ws.Range("V15").Value = 0.6363636363636364
$ws.Range("W15").Value = 0.561170113407071
$ws.Range("X15").Value = 0.0788800244827743
$ws.Range("Y15").Value = 7
$ws.Range("AA15").Value = 0.6550296044343662
$ws.Range("AB15").Value = 0.05511033388910008
$ws.Range("B16").Value = 0.002931070327758789
$ws.Range("C16").Value = 0.002305809691566049
$ws.Range("D16").Value = 0.0008890151977539063
$ws.Range("E16").Value = 0.000276430582129406
$ws.Range("R16").Value = 0.4912280701754386
$ws.Range("S16").Value = 0.4642857142857143
$ws.Range("T16").Value = 0.6060606060606061
$ws.Range("U16").Value = 0.5490196078431373
$ws.Range("V16").Value = 0.5937499999999999
$ws.Range("W16").Value = 0.5408687996729793
$ws.Range("X16").Value = 0.05557337426217009
$ws.Range("Y16").Value = 12
$ws.Range("AA16").Value = 0.6550296044343662
$ws.Range("AB16").Value = 0.05511033388910008
$ws.Range("B17").Value = 0.004493045806884766
$ws.Range("C17").Value = 0.005705922958679747
$ws.Range("D17").Value = 0.0007457733154296875
$ws.Range("E17").Value = 0.0001994385878857626
$ws.Range("R17").Value = 0.5
$ws.Range("S17").Value = 0.4363636363636363
$ws.Range("T17").Value = 0.6060606060606061
$ws.Range("U17").Value = 0.5599999999999999
$ws.Range("V17").Value = 0.5937499999999999
$ws.Range("W17").Value = 0.5392348484848485
$ws.Range("X17").Value = 0.06323074371958229
$ws.Range("Y17").Value = 13
$ws.Range("AA17").Value = 0.6550296044343662
$ws.Range("AB17").Value = 0.05511033388910008
